# Update parameter-estimate values on the model-fit sheets and the
# covariance-matrix values on the "* cov" sheets, per the new DTIC
# Falkson 1998 HR / survival-curve fit.

$wb = $excel.ActiveWorkbook

# --- weibull ---
$ws = $wb.Worksheets.Item("weibull")
$ws.Range("B2").Value = -2.84402157787172
$ws.Range("C2").Value = 0.166400404319523
$ws.Range("B3").Value = -0.0114217344621236
$ws.Range("C3").Value = 0.0918101730715185

# --- lognormal ---
$ws = $wb.Worksheets.Item("lognormal")
$ws.Range("B2").Value = 2.24237244057441
$ws.Range("C2").Value = 0.209726581516632
$ws.Range("B3").Value = -0.947676440017259
$ws.Range("C3").Value = 0.093789380646873

# --- llogis ---
$ws = $wb.Worksheets.Item("llogis")
$ws.Range("B2").Value = -2.33221811018217
$ws.Range("C2").Value = 0.13146167929006
$ws.Range("B3").Value = 0.487257647581758
$ws.Range("C3").Value = 0.110417509134753

# --- gompertz ---
$ws = $wb.Worksheets.Item("gompertz")
$ws.Range("B2").Value = -2.53774537520061
$ws.Range("C2").Value = 0.141652278771847
$ws.Range("B3").Value = -0.0256502815535124
$ws.Range("C3").Value = 0.0115962452167281

# --- exp --- (no data changes)

# --- weibull cov ---
$ws = $wb.Worksheets.Item("weibull cov")
$ws.Range("A2").Value = 0.0276890945577009
$ws.Range("B2").Value = -0.00936032305604423
$ws.Range("A3").Value = -0.00936032305604423
$ws.Range("B3").Value = 0.00842910787942218

# --- lognormal cov ---
$ws = $wb.Worksheets.Item("lognormal cov")
$ws.Range("A2").Value = 0.0439852389946526
$ws.Range("B2").Value = -0.0161834686181826
$ws.Range("A3").Value = -0.0161834686181826
$ws.Range("B3").Value = 0.00879644792212403

# --- llogis cov ---
$ws = $wb.Worksheets.Item("llogis cov")
$ws.Range("A2").Value = 0.0172821731217626
$ws.Range("B2").Value = 0.00414207079880718
$ws.Range("A3").Value = 0.00414207079880718
$ws.Range("B3").Value = 0.0121920263235233

# --- gompertz cov ---
$ws = $wb.Worksheets.Item("gompertz cov")
$ws.Range("A2").Value = 0.0200653680812571
$ws.Range("B2").Value = -0.000884315954813601
$ws.Range("A3").Value = -0.000884315954813601
$ws.Range("B3").Value = 0.00013447290312649

# --- exp cov --- (no data changes)

$wb.Save()
